$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet and update the header label text to reflect the new cutoff date (Oct 29 -> Oct 30)
$ws.Name = "Through 2021-10-30"
$ws.Range("B1").Value = "October 2021 (through October 30)"

# Helper to find a row number for a given neighborhood name in column A
function Get-RowForNeighborhood($name) {
    $cell = $ws.Columns.Item(1).Find($name, [Type]::Missing, [Type]::Missing, 1)
    return $cell.Row
}

# Map of neighborhood -> list of (column, newValue)
$changes = @{
    "Austin"              = @(@("BJ", 5))
    "North Lawndale"      = @(@("L", 12))
    "Little Italy, UIC"   = @(@("B", 2), @("AF", 1))
    "Humboldt Park"       = @(@("B", 11), @("L", 4))
    "Grand Crossing"      = @(@("AF", 2))
    "Roseland"            = @(@("AF", 2))
    "West Town"           = @(@("B", 9), @("L", 2))
    "Logan Square"        = @(@("V", 2))
    "Lower West Side"     = @(@("B", 4))
    "Grand Boulevard"     = @(@("AZ", 3))
    "West Pullman"        = @(@("B", 3))
    "Avondale"            = @(@("L", 1))
    "Kenwood"             = @(@("AP", 2))
    "Uptown"              = @(@("B", 2))
    "Washington Heights"  = @(@("L", 6))
    "Washington Park"     = @(@("AZ", 2))
    "South Chicago"       = @(@("AZ", 1))
    "Irving Park"         = @(@("L", 4))
    "New City"            = @(@("B", 5))
    "Albany Park"         = @(@("B", 2))
    "East Village"        = @(@("AF", 1))
    "Gage Park"           = @(@("V", 1))
    "Near South Side"     = @(@("L", 2))
    "South Deering"       = @(@("L", 3))
    "Wrigleyville"        = @(@("AF", 1))
}

foreach ($name in $changes.Keys) {
    $row = Get-RowForNeighborhood $name
    foreach ($pair in $changes[$name]) {
        $col = $pair[0]
        $val = $pair[1]
        $ws.Range("$col$row").Value = $val
    }
}
